$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds a "last updated" timestamp (serial date-time value) for each
# availability entry. This update re-runs the refresh job: the newest block
# of rows (2-15) gets a brand new timestamp, while the previous timestamps
# cascade down into the next two blocks (16-29 and 30-43), effectively
# shifting history down by one generation.

$newTimestamp = 44300.87771072562
$shiftedFromBlock1 = 44300.85556826389
$shiftedFromBlock2 = 44267.74495982639

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $newTimestamp
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $shiftedFromBlock1
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $shiftedFromBlock2
}
